# Insert a new data row at row 271 (shifting existing rows 271:367 down to 272:368)
# and populate it with the new Kiwi price record described in the commit.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 271; all rows below shift down by one,
# carrying their formatting (e.g. the date style on column D) along with them.
$ws.Rows.Item(271).Insert()

# Populate the newly inserted row 271 with the new record's values.
$ws.Cells.Item(271, 1).Value = 5
$ws.Cells.Item(271, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(271, 3).Value = "Maule"
$ws.Cells.Item(271, 4).Value = 44837
$ws.Cells.Item(271, 5).Value = 7
$ws.Cells.Item(271, 6).Value = "Fruta"
$ws.Cells.Item(271, 7).Value = 100101
$ws.Cells.Item(271, 8).Value = "Berries"
$ws.Cells.Item(271, 9).Value = 100101007
$ws.Cells.Item(271, 10).Value = "Kiwi"
$ws.Cells.Item(271, 11).Value = "Hayward"
$ws.Cells.Item(271, 12).Value = "Segunda"
$ws.Cells.Item(271, 13).Value = 280
$ws.Cells.Item(271, 14).Value = 8000
$ws.Cells.Item(271, 15).Value = 8000
$ws.Cells.Item(271, 16).Value = 8000
$ws.Cells.Item(271, 17).Value = "`$/bandeja 18 kilos"
$ws.Cells.Item(271, 18).Value = "Provincia de Curic$([char]0xF3)"
$ws.Cells.Item(271, 19).Value = 444
$ws.Cells.Item(271, 20).Value = 18
